$wb = $excel.ActiveWorkbook

# "Formula Samples" sheet (index 2 of 3: Constant Samples, Formula Samples, Test Case Samples)
$formula = $wb.Worksheets.Item("Formula Samples")

# Add a new "Sum" comparison row (row 11) that sums up the prior values via a
# formula, mirroring what the grading key does with SUM().
$formula.Range("A11").Value = "Sum"

# Update the comment in C4 to reflect the new wording that also mentions the
# "Sum" key comparison approach.
$formula.Range("C4").Value = "> This cell will be considered as wrong since ROUNDUP is a custom formula, and it multiplies against 0.01 instead of doing what the key does: dividing against 100. To compare this properly, use test cases. Or perhaps constant calculation will be enough."

$formula.Range("B11").Formula = "=B2+B3+B4"
$formula.Range("C11").Value = "> This cell will be considered as right, since the key uses SUM - basically expands the cell range and adds all of them. Algebra-wise, it will end up with similar result."

# "Test Case Samples" sheet loses the active-tab/selection state to the
# "Formula Samples" sheet above; record its new (non-active) selection first.
$testCase = $wb.Worksheets.Item("Test Case Samples")
[void]$testCase.Range("A12").Select()

# "Formula Samples" becomes the active/selected tab, with C11 selected.
[void]$formula.Activate()
[void]$formula.Range("C11").Select()
